# Applies per-row Price (D) and Volume(1h) (E) updates to match the refreshed
# crypto market snapshot described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.538.93"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "2.514.21"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").Value = "2.512.72"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("E10").Value = "  -2.27%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").Value = "2.973.92"
$ws.Range("E14").Value = "  -0.33%  "

$ws.Range("D15").Value = "69.337.28"
$ws.Range("E15").Value = "  -1.60%  "

$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").Value = "2.517.45"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.06%  "

$ws.Range("E22").Value = "  -1.52%  "

$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  -3.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "

$ws.Range("D28").Value = "2.644.17"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").Value = "0.0₃0888"
$ws.Range("E30").Value = "  -3.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "462.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.77%  "

$ws.Range("E33").Value = "  -5.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.30"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  -1.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.92%  "

$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("E45").Value = "  -7.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.577"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.69%  "
